$d = $word.ActiveDocument

$d.Content.Find.Execute("42+57=", $true, $true, $false, $false, $false, $true, 1, $false, "34+40=", 2) | Out-Null
$d.Content.Find.Execute("30+26=", $true, $true, $false, $false, $false, $true, 1, $false, "43-24=", 2) | Out-Null
$d.Content.Find.Execute("78-73=", $true, $true, $false, $false, $false, $true, 1, $false, "16+20=", 2) | Out-Null
$d.Content.Find.Execute("82-59=", $true, $true, $false, $false, $false, $true, 1, $false, "86-80=", 2) | Out-Null
$d.Content.Find.Execute("53-30=", $true, $true, $false, $false, $false, $true, 1, $false, "95-46=", 2) | Out-Null
$d.Content.Find.Execute("79-5=", $true, $true, $false, $false, $false, $true, 1, $false, "8+43=", 2) | Out-Null
$d.Content.Find.Execute("76-10=", $true, $true, $false, $false, $false, $true, 1, $false, "49+44=", 2) | Out-Null
$d.Content.Find.Execute("9+26=", $true, $true, $false, $false, $false, $true, 1, $false, "13+27=", 2) | Out-Null
$d.Content.Find.Execute("89-44=", $true, $true, $false, $false, $false, $true, 1, $false, "72+23=", 2) | Out-Null
$d.Content.Find.Execute("59-4=", $true, $true, $false, $false, $false, $true, 1, $false, "58+4=", 2) | Out-Null
$d.Content.Find.Execute("12-6=", $true, $true, $false, $false, $false, $true, 1, $false, "94+5=", 2) | Out-Null
$d.Content.Find.Execute("73+8=", $true, $true, $false, $false, $false, $true, 1, $false, "12+21=", 2) | Out-Null
$d.Content.Find.Execute("28+44=", $true, $true, $false, $false, $false, $true, 1, $false, "46-29=", 2) | Out-Null
$d.Content.Find.Execute("76-55=", $true, $true, $false, $false, $false, $true, 1, $false, "73-57=", 2) | Out-Null
$d.Content.Find.Execute("60-18=", $true, $true, $false, $false, $false, $true, 1, $false, "6+62=", 2) | Out-Null
$d.Content.Find.Execute("39-4=", $true, $true, $false, $false, $false, $true, 1, $false, "44+1=", 2) | Out-Null
$d.Content.Find.Execute("1+66=", $true, $true, $false, $false, $false, $true, 1, $false, "63-22=", 2) | Out-Null
$d.Content.Find.Execute("42-35=", $true, $true, $false, $false, $false, $true, 1, $false, "91-66=", 2) | Out-Null
$d.Content.Find.Execute("0+47=", $true, $true, $false, $false, $false, $true, 1, $false, "49-28=", 2) | Out-Null
$d.Content.Find.Execute("58+35=", $true, $true, $false, $false, $false, $true, 1, $false, "81-73=", 2) | Out-Null
$d.Content.Find.Execute("87-42=", $true, $true, $false, $false, $false, $true, 1, $false, "41+19=", 2) | Out-Null
$d.Content.Find.Execute("55+35=", $true, $true, $false, $false, $false, $true, 1, $false, "75+11=", 2) | Out-Null
$d.Content.Find.Execute("70+8=", $true, $true, $false, $false, $false, $true, 1, $false, "35+32=", 2) | Out-Null
$d.Content.Find.Execute("6+3=", $true, $true, $false, $false, $false, $true, 1, $false, "86-60=", 2) | Out-Null
$d.Content.Find.Execute("11+58=", $true, $true, $false, $false, $false, $true, 1, $false, "46+13=", 2) | Out-Null
$d.Content.Find.Execute("99-2=", $true, $true, $false, $false, $false, $true, 1, $false, "32+6=", 2) | Out-Null
$d.Content.Find.Execute("15-6=", $true, $true, $false, $false, $false, $true, 1, $false, "69-9=", 2) | Out-Null
$d.Content.Find.Execute("68-55=", $true, $true, $false, $false, $false, $true, 1, $false, "3+90=", 2) | Out-Null
$d.Content.Find.Execute("79-69=", $true, $true, $false, $false, $false, $true, 1, $false, "83-15=", 2) | Out-Null
$d.Content.Find.Execute("76-24=", $true, $true, $false, $false, $false, $true, 1, $false, "62+28=", 2) | Out-Null
$d.Content.Find.Execute("68-45=", $true, $true, $false, $false, $false, $true, 1, $false, "56-44=", 2) | Out-Null
$d.Content.Find.Execute("10+67=", $true, $true, $false, $false, $false, $true, 1, $false, "42-21=", 2) | Out-Null
$d.Content.Find.Execute("87-17=", $true, $true, $false, $false, $false, $true, 1, $false, "25+11=", 2) | Out-Null
$d.Content.Find.Execute("77+5=", $true, $true, $false, $false, $false, $true, 1, $false, "19-6=", 2) | Out-Null
$d.Content.Find.Execute("98-60=", $true, $true, $false, $false, $false, $true, 1, $false, "49+48=", 2) | Out-Null
$d.Content.Find.Execute("93-63=", $true, $true, $false, $false, $false, $true, 1, $false, "68+2=", 2) | Out-Null
$d.Content.Find.Execute("27+46=", $true, $true, $false, $false, $false, $true, 1, $false, "24+33=", 2) | Out-Null
$d.Content.Find.Execute("36+43=", $true, $true, $false, $false, $false, $true, 1, $false, "11+86=", 2) | Out-Null
$d.Content.Find.Execute("27+57=", $true, $true, $false, $false, $false, $true, 1, $false, "47-20=", 2) | Out-Null
$d.Content.Find.Execute("95-92=", $true, $true, $false, $false, $false, $true, 1, $false, "70+19=", 2) | Out-Null
$d.Content.Find.Execute("73-56=", $true, $true, $false, $false, $false, $true, 1, $false, "3+45=", 2) | Out-Null
$d.Content.Find.Execute("46-41=", $true, $true, $false, $false, $false, $true, 1, $false, "59-21=", 2) | Out-Null
$d.Content.Find.Execute("0+19=", $true, $true, $false, $false, $false, $true, 1, $false, "64-51=", 2) | Out-Null
$d.Content.Find.Execute("53-39=", $true, $true, $false, $false, $false, $true, 1, $false, "64-1=", 2) | Out-Null
$d.Content.Find.Execute("55+30=", $true, $true, $false, $false, $false, $true, 1, $false, "44+3=", 2) | Out-Null
$d.Content.Find.Execute("40+34=", $true, $true, $false, $false, $false, $true, 1, $false, "34+16=", 2) | Out-Null
$d.Content.Find.Execute("85-39=", $true, $true, $false, $false, $false, $true, 1, $false, "48-37=", 2) | Out-Null
$d.Content.Find.Execute("57-52=", $true, $true, $false, $false, $false, $true, 1, $false, "6+16=", 2) | Out-Null
$d.Content.Find.Execute("42+0=", $true, $true, $false, $false, $false, $true, 1, $false, "95-35=", 2) | Out-Null
$d.Content.Find.Execute("21+76=", $true, $true, $false, $false, $false, $true, 1, $false, "83-7=", 2) | Out-Null
$d.Content.Find.Execute("42-19=", $true, $true, $false, $false, $false, $true, 1, $false, "28+65=", 2) | Out-Null
$d.Content.Find.Execute("80-33=", $true, $true, $false, $false, $false, $true, 1, $false, "84-76=", 2) | Out-Null
$d.Content.Find.Execute("91-58=", $true, $true, $false, $false, $false, $true, 1, $false, "11+62=", 2) | Out-Null
$d.Content.Find.Execute("74-33=", $true, $true, $false, $false, $false, $true, 1, $false, "78-59=", 2) | Out-Null
$d.Content.Find.Execute("72-54=", $true, $true, $false, $false, $false, $true, 1, $false, "42+44=", 2) | Out-Null
$d.Content.Find.Execute("67+10=", $true, $true, $false, $false, $false, $true, 1, $false, "76-20=", 2) | Out-Null
$d.Content.Find.Execute("70-56=", $true, $true, $false, $false, $false, $true, 1, $false, "61+35=", 2) | Out-Null
$d.Content.Find.Execute("53-34=", $true, $true, $false, $false, $false, $true, 1, $false, "76-62=", 2) | Out-Null
$d.Content.Find.Execute("81-34=", $true, $true, $false, $false, $false, $true, 1, $false, "22+76=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $true, $true, $false, $false, $false, $true, 1, $false, "47-38=", 2) | Out-Null
$d.Content.Find.Execute("27+58=", $true, $true, $false, $false, $false, $true, 1, $false, "26+8=", 2) | Out-Null
$d.Content.Find.Execute("14-6=", $true, $true, $false, $false, $false, $true, 1, $false, "55-14=", 2) | Out-Null
$d.Content.Find.Execute("82-51=", $true, $true, $false, $false, $false, $true, 1, $false, "5+79=", 2) | Out-Null
$d.Content.Find.Execute("75-18=", $true, $true, $false, $false, $false, $true, 1, $false, "11+27=", 2) | Out-Null
$d.Content.Find.Execute("94-8=", $true, $true, $false, $false, $false, $true, 1, $false, "10+0=", 2) | Out-Null
$d.Content.Find.Execute("83-16=", $true, $true, $false, $false, $false, $true, 1, $false, "22+67=", 2) | Out-Null
$d.Content.Find.Execute("88-65=", $true, $true, $false, $false, $false, $true, 1, $false, "52-8=", 2) | Out-Null
$d.Content.Find.Execute("18+42=", $true, $true, $false, $false, $false, $true, 1, $false, "61+11=", 2) | Out-Null
$d.Content.Find.Execute("55+2=", $true, $true, $false, $false, $false, $true, 1, $false, "54+7=", 2) | Out-Null
$d.Content.Find.Execute("9+52=", $true, $true, $false, $false, $false, $true, 1, $false, "42+26=", 2) | Out-Null
$d.Content.Find.Execute("52+11=", $true, $true, $false, $false, $false, $true, 1, $false, "31+60=", 2) | Out-Null
$d.Content.Find.Execute("12+87=", $true, $true, $false, $false, $false, $true, 1, $false, "93-13=", 2) | Out-Null
$d.Content.Find.Execute("85-63=", $true, $true, $false, $false, $false, $true, 1, $false, "25-18=", 2) | Out-Null
$d.Content.Find.Execute("73-72=", $true, $true, $false, $false, $false, $true, 1, $false, "49-37=", 2) | Out-Null
$d.Content.Find.Execute("25+1=", $true, $true, $false, $false, $false, $true, 1, $false, "42-17=", 2) | Out-Null
$d.Content.Find.Execute("25-17=", $true, $true, $false, $false, $false, $true, 1, $false, "74+4=", 2) | Out-Null
$d.Content.Find.Execute("54-38=", $true, $true, $false, $false, $false, $true, 1, $false, "78-24=", 2) | Out-Null
$d.Content.Find.Execute("72-59=", $true, $true, $false, $false, $false, $true, 1, $false, "28+17=", 2) | Out-Null
$d.Content.Find.Execute("58-49=", $true, $true, $false, $false, $false, $true, 1, $false, "69+8=", 2) | Out-Null
$d.Content.Find.Execute("50+20=", $true, $true, $false, $false, $false, $true, 1, $false, "6+23=", 2) | Out-Null
$d.Content.Find.Execute("55-38=", $true, $true, $false, $false, $false, $true, 1, $false, "59+8=", 2) | Out-Null
$d.Content.Find.Execute("6-5=", $true, $true, $false, $false, $false, $true, 1, $false, "56+43=", 2) | Out-Null
$d.Content.Find.Execute("19+74=", $true, $true, $false, $false, $false, $true, 1, $false, "1+79=", 2) | Out-Null
$d.Content.Find.Execute("77-49=", $true, $true, $false, $false, $false, $true, 1, $false, "13+42=", 2) | Out-Null
$d.Content.Find.Execute("33+1=", $true, $true, $false, $false, $false, $true, 1, $false, "77+13=", 2) | Out-Null
$d.Content.Find.Execute("93-40=", $true, $true, $false, $false, $false, $true, 1, $false, "99-5=", 2) | Out-Null
$d.Content.Find.Execute("88+4=", $true, $true, $false, $false, $false, $true, 1, $false, "29+52=", 2) | Out-Null
$d.Content.Find.Execute("73-42=", $true, $true, $false, $false, $false, $true, 1, $false, "93-16=", 2) | Out-Null
$d.Content.Find.Execute("47-24=", $true, $true, $false, $false, $false, $true, 1, $false, "94-87=", 2) | Out-Null
$d.Content.Find.Execute("79-61=", $true, $true, $false, $false, $false, $true, 1, $false, "53+25=", 2) | Out-Null
$d.Content.Find.Execute("28+22=", $true, $true, $false, $false, $false, $true, 1, $false, "31-1=", 2) | Out-Null
$d.Content.Find.Execute("69-47=", $true, $true, $false, $false, $false, $true, 1, $false, "48-20=", 2) | Out-Null
$d.Content.Find.Execute("2+11=", $true, $true, $false, $false, $false, $true, 1, $false, "53+14=", 2) | Out-Null
$d.Content.Find.Execute("18+70=", $true, $true, $false, $false, $false, $true, 1, $false, "34-33=", 2) | Out-Null
$d.Content.Find.Execute("48-39=", $true, $true, $false, $false, $false, $true, 1, $false, "29+19=", 2) | Out-Null
$d.Content.Find.Execute("47-23=", $true, $true, $false, $false, $false, $true, 1, $false, "10+73=", 2) | Out-Null
$d.Content.Find.Execute("44+7=", $true, $true, $false, $false, $false, $true, 1, $false, "57-19=", 2) | Out-Null
$d.Content.Find.Execute("21+33=", $true, $true, $false, $false, $false, $true, 1, $false, "39+37=", 2) | Out-Null
$d.Content.Find.Execute("94-22=", $true, $true, $false, $false, $false, $true, 1, $false, "52-1=", 2) | Out-Null
$d.Content.Find.Execute("73-1=", $true, $true, $false, $false, $false, $true, 1, $false, "47-22=", 2) | Out-Null
